# Update "想去人数" (interested count, column F) values on the "展览"
# and "全部类型" sheets, matching the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 368
$wsExpo.Range("F3").Value = 69
$wsExpo.Range("F4").Value = 277
$wsExpo.Range("F5").Value = 4116
$wsExpo.Range("F6").Value = 39
$wsExpo.Range("F7").Value = 453

# --- Sheet "全部类型" (all types, combined list) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 368
$wsAll.Range("F3").Value = 69
$wsAll.Range("F4").Value = 277
$wsAll.Range("F5").Value = 4116
$wsAll.Range("F8").Value = 39
$wsAll.Range("F9").Value = 453
